$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sheet / workbook title from "12-03" to "12-04"
$ws.Name = "Through 2021-12-04"

# Update the row label for December
$ws.Range("A13").Value = "December (through 12-04)"

# Update the December row (row 13) values for each year column (B..H)
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 13
$ws.Range("D13").Value = 16
$ws.Range("E13").Value = 11
$ws.Range("F13").Value = 4
$ws.Range("G13").Value = 19
$ws.Range("H13").Value = 31

# Update the Total row (row 14) values for each year column (B..H)
$ws.Range("B14").Value = 293
$ws.Range("C14").Value = 576
$ws.Range("D14").Value = 837
$ws.Range("E14").Value = 693
$ws.Range("F14").Value = 538
$ws.Range("G14").Value = 1283
$ws.Range("H14").Value = 1675
